$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "749÷6="
$t.Cell(1, 2).Range.Text = "115÷9="
$t.Cell(1, 3).Range.Text = "673÷4="
$t.Cell(1, 4).Range.Text = "847÷9="
$t.Cell(1, 5).Range.Text = "702÷6="
$t.Cell(5, 1).Range.Text = "176÷8="
$t.Cell(5, 2).Range.Text = "342÷6="
$t.Cell(5, 3).Range.Text = "766÷7="
$t.Cell(5, 4).Range.Text = "668÷3="
$t.Cell(5, 5).Range.Text = "270÷2="
$t.Cell(9, 1).Range.Text = "973÷7="
$t.Cell(9, 2).Range.Text = "926÷8="
$t.Cell(9, 3).Range.Text = "975÷4="
$t.Cell(9, 4).Range.Text = "728÷5="
$t.Cell(9, 5).Range.Text = "786÷5="
$t.Cell(13, 1).Range.Text = "952÷9="
$t.Cell(13, 2).Range.Text = "817÷5="
$t.Cell(13, 3).Range.Text = "229÷6="
$t.Cell(13, 4).Range.Text = "802÷6="
$t.Cell(13, 5).Range.Text = "466÷3="
$t.Cell(17, 1).Range.Text = "839÷2="
$t.Cell(17, 2).Range.Text = "935÷4="
$t.Cell(17, 3).Range.Text = "729÷8="
$t.Cell(17, 4).Range.Text = "788÷4="
$t.Cell(17, 5).Range.Text = "565÷2="
